$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 448.15384
$ws.Range("I9").Value = 440.25
$ws.Range("K9").Value = 440.25
$ws.Range("M9").Value = -271.25

$ws.Range("H15").Value = 1073.6154
$ws.Range("I15").Value = 1073.6154
$ws.Range("K15").Value = 3220.8462
$ws.Range("M15").Value = -3051.8462

$ws.Range("H51").Value = 4470741.5
$ws.Range("I51").Value = 3849.75
$ws.Range("J51").Value = 5959705.5
$ws.Range("K51").Value = 3849.75
$ws.Range("L51").Value = 5959705.5
$ws.Range("M51").Value = -3365.75
$ws.Range("N51").Value = -5960673.5

$ws.Range("H58").Value = 1137.5
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300

$ws.Range("H88").Value = 1118.8695
$ws.Range("I88").Value = 1149.1666
$ws.Range("J88").Value = 1108.1765
$ws.Range("K88").Value = 1149.1666
$ws.Range("L88").Value = 1108.1765
$ws.Range("M88").Value = -743.1666
$ws.Range("N88").Value = -1920.1765

$ws.Range("H91").Value = 1118.8695
$ws.Range("I91").Value = 1149.1666
$ws.Range("J91").Value = 1108.1765
$ws.Range("K91").Value = 1149.1666
$ws.Range("L91").Value = 1108.1765
$ws.Range("M91").Value = 254.8334
$ws.Range("N91").Value = -3916.1765

$ws.Range("H113").Value = 83336240
$ws.Range("I113").Value = 500001000
$ws.Range("J113").Value = 3292.8
$ws.Range("K113").Value = 500001000
$ws.Range("L113").Value = 3292.8
$ws.Range("M113").Value = -499997746
$ws.Range("N113").Value = -9800.799999999999

$ws.Range("H132").Value = 2111.0605
$ws.Range("I132").Value = 1907.2572
$ws.Range("J132").Value = 2603
$ws.Range("K132").Value = 5721.7716
$ws.Range("L132").Value = 7809
$ws.Range("M132").Value = -3191.7716
$ws.Range("N132").Value = -12869

$ws.Range("H138").Value = 2041.2653
$ws.Range("I138").Value = 1455.9117
$ws.Range("J138").Value = 3368.0667
$ws.Range("K138").Value = 4367.7351
$ws.Range("L138").Value = 10104.2001
$ws.Range("M138").Value = 772.2649000000001
$ws.Range("N138").Value = -20384.2001

$ws.Range("H141").Value = 2212.3215
$ws.Range("I141").Value = 2151.7693
$ws.Range("K141").Value = 6455.3079
$ws.Range("M141").Value = -1275.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 12622.223
$ws.Range("I6").Value = 13957.143
$ws.Range("J6").Value = 7950
$ws.Range("K6").Value = 13957.143
$ws.Range("L6").Value = 7950
$ws.Range("M6").Value = -13784.143
$ws.Range("N6").Value = -8296

$ws.Range("H26").Value = 5374.75
$ws.Range("I26").Value = 5374.75
$ws.Range("K26").Value = 5374.75
$ws.Range("M26").Value = -5044.75

$ws.Range("H32").Value = 129816.164
$ws.Range("I32").Value = 144260.75
$ws.Range("K32").Value = 144260.75
$ws.Range("M32").Value = -143973.75

$ws.Range("H38").Value = 20000
$ws.Range("I38").Value = 20000
$ws.Range("K38").Value = 20000
$ws.Range("M38").Value = -19533

$ws.Range("H39").Value = 6399
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 6399
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 6399
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -7439

$ws.Range("H61").Value = 1055401.6
$ws.Range("I61").Value = 910693.75
$ws.Range("K61").Value = 910693.75
$ws.Range("M61").Value = -910481.75

$ws.Range("H74").Value = 2969.7058
$ws.Range("I74").Value = 2940.6
$ws.Range("K74").Value = 2940.6
$ws.Range("M74").Value = -2066.6

$ws.Range("H77").Value = 2969.7058
$ws.Range("I77").Value = 2940.6
$ws.Range("K77").Value = 14703
$ws.Range("M77").Value = -10335

$ws.Range("H110").Value = 41667670
$ws.Range("I110").Value = 45455504
$ws.Range("K110").Value = 45455504
$ws.Range("M110").Value = -45453459

$ws.Range("H132").Value = 589797.5
$ws.Range("I132").Value = 358564.8
$ws.Range("K132").Value = 1075694.4
$ws.Range("M132").Value = -1073164.4

$ws.Range("H136").Value = 1055401.6
$ws.Range("I136").Value = 910693.75
$ws.Range("K136").Value = 2732081.25
$ws.Range("M136").Value = -2729531.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1649.0605
$ws.Range("J134").Value = 2479.6
$ws.Range("L134").Value = 7438.799999999999
$ws.Range("N134").Value = -12508.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2838.0122
$ws.Range("I31").Value = 1990.2727
$ws.Range("J31").Value = 2969.352
$ws.Range("K31").Value = 1990.2727
$ws.Range("L31").Value = 2969.352
$ws.Range("M31").Value = -1695.2727
$ws.Range("N31").Value = -3559.352

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H34").Value = 2838.0122
$ws.Range("I34").Value = 1990.2727
$ws.Range("J34").Value = 2969.352
$ws.Range("K34").Value = 1990.2727
$ws.Range("L34").Value = 2969.352
$ws.Range("M34").Value = -1788.2727
$ws.Range("N34").Value = -3373.352

$ws.Range("H86").Value = 4057.0625
$ws.Range("I86").Value = 4269.909
$ws.Range("K86").Value = 4269.909
$ws.Range("M86").Value = -3146.909

$ws.Range("H89").Value = 4057.0625
$ws.Range("I89").Value = 4269.909
$ws.Range("K89").Value = 21349.545
$ws.Range("M89").Value = -15733.545

$ws.Range("H98").Value = 64995
$ws.Range("I98").Value = 64995
$ws.Range("K98").Value = 64995
$ws.Range("M98").Value = -62749

$ws.Range("H132").Value = 1665.9333
$ws.Range("I132").Value = 1378.6207
$ws.Range("K132").Value = 4135.8621
$ws.Range("M132").Value = -1605.8621

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1282.4445
$ws.Range("I23").Value = 621.3333
$ws.Range("K23").Value = 1863.9999
$ws.Range("M23").Value = -1628.9999

$ws.Range("H80").Value = 2749.5
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 2749.5
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3208.25
$ws.Range("I80").Value = 3075.6667
$ws.Range("J80").Value = 3340.8333
$ws.Range("K80").Value = 3075.6667
$ws.Range("L80").Value = 3340.8333
$ws.Range("M80").Value = -2077.6667
$ws.Range("N80").Value = -5336.8333

$ws.Range("H83").Value = 3208.25
$ws.Range("I83").Value = 3075.6667
$ws.Range("J83").Value = 3340.8333
$ws.Range("K83").Value = 15378.3335
$ws.Range("L83").Value = 16704.1665
$ws.Range("M83").Value = -10386.3335
$ws.Range("N83").Value = -26688.1665

$ws.Range("H102").Value = 3286.2
$ws.Range("I102").Value = 2143.6667
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 2143.6667
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -521.6667000000002
$ws.Range("N102").Value = -8244

$ws.Range("H126").Value = 8396.286
$ws.Range("I126").Value = 11782.546
$ws.Range("K126").Value = 35347.638
$ws.Range("M126").Value = -32877.638

$ws.Range("H132").Value = 253288.5
$ws.Range("I132").Value = 305803.44
$ws.Range("K132").Value = 917410.3200000001
$ws.Range("M132").Value = -914880.3200000001

$ws.Range("H139").Value = 72500
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 72500
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 72500
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -82780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16132151
$ws.Range("I7").Value = 26318472
$ws.Range("K7").Value = 26318472
$ws.Range("M7").Value = -26318360

$ws.Range("H32").Value = 3113
$ws.Range("I32").Value = 3113
$ws.Range("K32").Value = 3113
$ws.Range("M32").Value = -2796

$ws.Range("H40").Value = 3395
$ws.Range("I40").Value = 3015.2307
$ws.Range("J40").Value = 4012.125
$ws.Range("K40").Value = 3015.2307
$ws.Range("L40").Value = 4012.125
$ws.Range("M40").Value = -2879.2307
$ws.Range("N40").Value = -4284.125

$ws.Range("H122").Value = 4799.4546
$ws.Range("I122").Value = 4383.3335
$ws.Range("J122").Value = 5298.8
$ws.Range("K122").Value = 13150.0005
$ws.Range("L122").Value = 15896.4
$ws.Range("M122").Value = -10700.0005
$ws.Range("N122").Value = -20796.4

$ws.Range("H126").Value = 16132151
$ws.Range("I126").Value = 26318472
$ws.Range("K126").Value = 78955416
$ws.Range("M126").Value = -78952946

$ws.Range("H136").Value = 6211.8066
$ws.Range("I136").Value = 3214.15
$ws.Range("K136").Value = 9642.450000000001
$ws.Range("M136").Value = -7092.450000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4400.533
$ws.Range("I96").Value = 3650
$ws.Range("J96").Value = 5526.3335
$ws.Range("K96").Value = 3650
$ws.Range("L96").Value = 5526.3335
$ws.Range("M96").Value = -2277
$ws.Range("N96").Value = -8272.333500000001

$ws.Range("H107").Value = 495.1111
$ws.Range("I107").Value = 500.30768
$ws.Range("K107").Value = 1500.92304
$ws.Range("M107").Value = 419.0769599999999

$ws.Range("H122").Value = 2677.8
$ws.Range("I122").Value = 1849
$ws.Range("J122").Value = 3625
$ws.Range("K122").Value = 5547
$ws.Range("L122").Value = 10875
$ws.Range("M122").Value = -3097
$ws.Range("N122").Value = -15775

$ws.Range("H132").Value = 543410.5
$ws.Range("I132").Value = 590968.8
$ws.Range("J132").Value = 4416.3335
$ws.Range("K132").Value = 1772906.4
$ws.Range("L132").Value = 13249.0005
$ws.Range("M132").Value = -1770376.4
$ws.Range("N132").Value = -18309.0005

$ws.Range("H136").Value = 3051.0613
$ws.Range("I136").Value = 2454.4722
$ws.Range("J136").Value = 4703.154
$ws.Range("K136").Value = 7363.4166
$ws.Range("L136").Value = 14109.462
$ws.Range("M136").Value = -4813.4166
$ws.Range("N136").Value = -19209.462
